$wb = $excel.ActiveWorkbook

# --- Sheet "LoginTestData": update the password test value ---
$ws1 = $wb.Worksheets.Item("LoginTestData")
$ws1.Activate()
$ws1.Range("B1").Value = "Nada_1234"
[void]$ws1.Range("B2").Select()

# --- Sheet "RegisterTestData": update the registration email test value ---
$ws2 = $wb.Worksheets.Item("RegisterTestData")
$ws2.Activate()
$ws2.Range("A1").Value = "nadasalama4@gmail.com"
[void]$ws2.Range("B3").Select()
